$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (rows 2 through 267) holds a date value (serial 45182) that
# must be updated to serial 45184 for every data row.
$lastRow = $ws.Range("A1").End(4).Row
if ($lastRow -lt 267) { $lastRow = 267 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45182) {
        $cell.Value2 = 45184
    }
}
